$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 7.240679918064295
$ws.Range("D2").Value = 5.821966626816713
$ws.Range("E2").Value = 12.46065815797079
$ws.Range("F2").Value = 53.88144777410825
$ws.Range("G2").Value = 71.55053117575936
$ws.Range("H2").Value = 24.17023702721602
$ws.Range("I2").Value = 40.9835502479899
$ws.Range("J2").Value = 10.92775514833046
$ws.Range("L2").Value = 10.45052131975529

$ws.Range("C3").Value = 7.253979846252999
$ws.Range("D3").Value = 5.803469237481298
$ws.Range("E3").Value = 12.48958818472225
$ws.Range("F3").Value = 53.47938016116055
$ws.Range("G3").Value = 70.64340532633257
$ws.Range("H3").Value = 24.08215247237879
$ws.Range("I3").Value = 40.73299499319318
$ws.Range("J3").Value = 10.9529646677571
$ws.Range("L3").Value = 10.47620201549809

$ws.Range("C4").Value = 7.26271621327336
$ws.Range("D4").Value = 5.791878923515848
$ws.Range("E4").Value = 12.50899615333631
$ws.Range("F4").Value = 53.24722460944307
$ws.Range("G4").Value = 70.10457091085938
$ws.Range("H4").Value = 24.03432458196458
$ws.Range("I4").Value = 40.59028273905891
$ws.Range("J4").Value = 10.97000095944364
$ws.Range("L4").Value = 10.49323879232863

$ws.Range("C5").Value = 7.266420037099552
$ws.Range("D5").Value = 5.78709693883918
$ws.Range("E5").Value = 12.51731877923824
$ws.Range("F5").Value = 53.15639101174025
$ws.Range("G5").Value = 69.88979535565961
$ws.Range("H5").Value = 24.01641548997419
$ws.Range("I5").Value = 40.53496255053431
$ws.Range("J5").Value = 10.97733491058786
$ws.Range("L5").Value = 10.50050077123841

$ws.Range("C6").Value = 7.267043741598628
$ws.Range("D6").Value = 5.786299331798079
$ws.Range("E6").Value = 12.51872573618098
$ws.Range("F6").Value = 53.14153789802576
$ws.Range("G6").Value = 69.85442899561012
$ws.Range("H6").Value = 24.0135373898203
$ws.Range("I6").Value = 40.52594897505884
$ws.Range("J6").Value = 10.9785763507118
$ws.Range("L6").Value = 10.50172591506588

$ws.Range("C7").Value = 7.262765582146261
$ws.Range("D7").Value = 5.791814670616978
$ws.Range("E7").Value = 12.50910671986533
$ws.Range("F7").Value = 53.24598423564467
$ws.Range("G7").Value = 70.10165461186159
$ws.Range("H7").Value = 24.03407664256475
$ws.Range("I7").Value = 40.5895251419211
$ws.Range("J7").Value = 10.97009828257198
$ws.Range("L7").Value = 10.49333543642013

$ws.Range("C8").Value = 7.245147584407692
$ws.Range("D8").Value = 5.815636878004721
$ws.Range("E8").Value = 12.47029192098302
$ws.Range("F8").Value = 53.73980567169957
$ws.Range("G8").Value = 71.23413647710709
$ws.Range("H8").Value = 24.13857231598908
$ws.Range("I8").Value = 40.89487285550843
$ws.Range("J8").Value = 10.93612404254448
$ws.Range("L8").Value = 10.45911289237452

$ws.Range("C9").Value = 7.215108822381551
$ws.Range("D9").Value = 5.860518890367795
$ws.Range("E9").Value = 12.40722232928423
$ws.Range("F9").Value = 54.82147432476124
$ws.Range("G9").Value = 73.58779030816869
$ws.Range("H9").Value = 24.3926933820049
$ws.Range("I9").Value = 41.58009133903971
$ws.Range("J9").Value = 10.88186553487734
$ws.Range("L9").Value = 10.40205533037947

$ws.Range("C10").Value = 7.195769984867233
$ws.Range("D10").Value = 5.892397699004501
$ws.Range("E10").Value = 12.36883266751558
$ws.Range("F10").Value = 55.68014242145668
$ws.Range("G10").Value = 75.382615947738
$ws.Range("H10").Value = 24.60868663556985
$ws.Range("I10").Value = 42.1333977960884
$ws.Range("J10").Value = 10.84954933826324
$ws.Range("L10").Value = 10.36624493908684

$ws.Range("C11").Value = 7.187561265708159
$ws.Range("D11").Value = 5.90666691072355
$ws.Range("E11").Value = 12.35309286705555
$ws.Range("F11").Value = 56.08346203390886
$ws.Range("G11").Value = 76.21001617407882
$ws.Range("H11").Value = 24.71311884555085
$ws.Range("I11").Value = 42.39525773261945
$ws.Range("J11").Value = 10.83648883126477
$ws.Range("L11").Value = 10.3512767649242

$ws.Range("C12").Value = 7.184537182689876
$ws.Range("D12").Value = 5.912037162253357
$ws.Range("E12").Value = 12.34738042714578
$ws.Range("F12").Value = 56.23790833005967
$ws.Range("G12").Value = 76.52462051688899
$ws.Range("H12").Value = 24.75353391884828
$ws.Range("I12").Value = 42.49581509440115
$ws.Range("J12").Value = 10.83177920044877
$ws.Range("L12").Value = 10.34579855375384

$ws.Range("C13").Value = 7.185184724326795
$ws.Range("D13").Value = 5.910882054456971
$ws.Range("E13").Value = 12.3485996788087
$ws.Range("F13").Value = 56.20457092824645
$ws.Range("G13").Value = 76.45681216362927
$ws.Range("H13").Value = 24.74479148950042
$ws.Range("I13").Value = 42.47409719890808
$ws.Range("J13").Value = 10.83278299952347
$ws.Range("L13").Value = 10.34696994147893

$ws.Range("C14").Value = 7.187310782939234
$ws.Range("D14").Value = 5.907109388213586
$ws.Range("E14").Value = 12.35261793331725
$ws.Range("F14").Value = 56.0961346323373
$ws.Range("G14").Value = 76.23587471576265
$ws.Range("H14").Value = 24.71642652038837
$ws.Range("I14").Value = 42.40350303253044
$ws.Range("J14").Value = 10.83609663402406
$ws.Range("L14").Value = 10.35082226427562

$ws.Range("C15").Value = 7.188624036432528
$ws.Range("D15").Value = 5.904794206737927
$ws.Range("E15").Value = 12.35511151146542
$ws.Range("F15").Value = 56.02993468999211
$ws.Range("G15").Value = 76.10070318765429
$ws.Range("H15").Value = 24.69916471176806
$ws.Range("I15").Value = 42.36044200767272
$ws.Range("J15").Value = 10.83815708467471
$ws.Range("L15").Value = 10.35320664912932

$ws.Range("C16").Value = 7.196318265947841
$ws.Range("D16").Value = 5.891460513898409
$ws.Range("E16").Value = 12.3698959881672
$ws.Range("F16").Value = 55.65403096743179
$ws.Range("G16").Value = 75.32873880181516
$ws.Range("H16").Value = 24.60198441261328
$ws.Range("I16").Value = 42.11648363326403
$ws.Range("J16").Value = 10.8504358995421
$ws.Range("L16").Value = 10.3672497312335

$ws.Range("C17").Value = 7.201188997336139
$ws.Range("D17").Value = 5.883221331877565
$ws.Range("E17").Value = 12.37940729992341
$ws.Range("F17").Value = 55.42660414023241
$ws.Range("G17").Value = 74.85775819926107
$ws.Range("H17").Value = 24.54393606585003
$ws.Range("I17").Value = 41.96938018930092
$ws.Range("J17").Value = 10.85838883625912
$ws.Range("L17").Value = 10.37620318873229

$ws.Range("C18").Value = 7.204045931745183
$ws.Range("D18").Value = 5.878460461598061
$ws.Range("E18").Value = 12.38504021210547
$ws.Range("F18").Value = 55.29699710012249
$ws.Range("G18").Value = 74.58790791313933
$ws.Range("H18").Value = 24.5111307440996
$ws.Range("I18").Value = 41.8857303390884
$ws.Range("J18").Value = 10.86311753393247
$ws.Range("L18").Value = 10.38147744019919

$ws.Range("C19").Value = 7.205022766883536
$ws.Range("D19").Value = 5.876844741607141
$ws.Range("E19").Value = 12.3869752874138
$ws.Range("F19").Value = 55.25332424940258
$ws.Range("G19").Value = 74.49672950215754
$ws.Range("H19").Value = 24.50012403820394
$ws.Range("I19").Value = 41.85757474748246
$ws.Range("J19").Value = 10.86474509740397
$ws.Range("L19").Value = 10.3832845937377

$ws.Range("C20").Value = 7.200664765968572
$ws.Range("D20").Value = 5.884100672650344
$ws.Range("E20").Value = 12.37837801159862
$ws.Range("F20").Value = 55.45069040463446
$ws.Range("G20").Value = 74.90778877801739
$ws.Range("H20").Value = 24.55005525359931
$ws.Range("I20").Value = 41.98494068367553
$ws.Range("J20").Value = 10.85752625329142
$ws.Range("L20").Value = 10.37523719837469

$ws.Range("C21").Value = 7.186684019975939
$ws.Range("D21").Value = 5.908218410409638
$ws.Range("E21").Value = 12.35143094719349
$ws.Range("F21").Value = 56.12793926527944
$ws.Range("G21").Value = 76.30073676084916
$ws.Range("H21").Value = 24.72473457970336
$ws.Range("I21").Value = 42.42420086915503
$ws.Range("J21").Value = 10.83511692977272
$ws.Range("L21").Value = 10.34968559119878

$ws.Range("C22").Value = 7.178038512520907
$ws.Range("D22").Value = 5.92378741528628
$ws.Range("E22").Value = 12.33526432029696
$ws.Range("F22").Value = 56.58052096965572
$ws.Range("G22").Value = 77.21847383996257
$ws.Range("H22").Value = 24.84395200054333
$ws.Range("I22").Value = 42.71939213923248
$ws.Range("J22").Value = 10.82184741418479
$ws.Range("L22").Value = 10.33409299526906

$ws.Range("C23").Value = 7.18260787526857
$ws.Range("D23").Value = 5.915495547875958
$ws.Range("E23").Value = 12.34376055245099
$ws.Range("F23").Value = 56.33809516110065
$ws.Range("G23").Value = 76.72807974620704
$ws.Range("H23").Value = 24.77986769463016
$ws.Range("I23").Value = 42.56112286494065
$ws.Range("J23").Value = 10.82880360110849
$ws.Range("L23").Value = 10.34231384795046

$ws.Range("C24").Value = 7.20090159438558
$ws.Range("D24").Value = 5.883703197632067
$ws.Range("E24").Value = 12.37884283974535
$ws.Range("F24").Value = 55.4397974339599
$ws.Range("G24").Value = 74.88516705758192
$ws.Range("H24").Value = 24.54728699991248
$ws.Range("I24").Value = 41.97790290789926
$ws.Range("J24").Value = 10.85791573976633
$ws.Range("L24").Value = 10.37567352767844

$ws.Range("C25").Value = 7.222754232826439
$ws.Range("D25").Value = 5.848571135465708
$ws.Range("E25").Value = 12.42288824186719
$ws.Range("F25").Value = 54.51722893744711
$ws.Range("G25").Value = 72.93841039704155
$ws.Range("H25").Value = 24.31874453786628
$ws.Range("I25").Value = 41.38574222493581
$ws.Range("J25").Value = 10.89521915254888
$ws.Range("L25").Value = 10.41641664866664
